$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update existing cells whose text changes (B2, A3, B3)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "password"
$ws.Range("A3").Value = "rahulscreencast9892@gmail.com"
$ws.Range("B3").Value = "Ra987456321@"

# ---------------------------------------------------------------------------
# 2. Fill in the brand-new header row (row 2, columns C:H) and data row
#    (row 3, columns C:H) with the additional test-data columns.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "url"
$ws.Range("D2").Value = "browser"
$ws.Range("E2").Value = "firstName"
$ws.Range("F2").Value = "lastName"
$ws.Range("G2").Value = "email"
$ws.Range("H2").Value = "description"

$ws.Range("C3").Value = "https://ui.cogmento.com/"
$ws.Range("D3").Value = "chrome"
$ws.Range("E3").Value = "Sumeet"
$ws.Range("F3").Value = "Desai"
$ws.Range("G3").Value = "sumeet.desai@gmail.com"
$ws.Range("H3").Value = "Create a follow up activity"

# ---------------------------------------------------------------------------
# 3. Normal (bordered, non-hyperlink) formatting for the new plain cells --
#    copy the format already used across row 2 / row 3 (cell B2 carries the
#    plain bordered style) onto every newly written cell that isn't a
#    hyperlink.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("C2:H2").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3:F3").PasteSpecial(-4122)
$ws.Range("H3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Hyperlinks. The engine's Range.Hyperlinks.Delete() clears every
#    hyperlink on the sheet, so drop them all and recreate the full set in
#    document order (this also reproduces the rId1..rId5 ordering seen in
#    the target file).
#
#    Hyperlinks.Add() always stamps a brand-new cell style, so first give
#    every hyperlink target cell the *same* plain bordered style (copied
#    from B2); that way all five Add() calls share/re-use one throw-away
#    style instead of minting five, keeping styles.xml tidy. The real
#    bordered-hyperlink look is restored by the PasteSpecial pass below.
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("B2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:rahulscreencast9892@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:test@rahul.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://ui.cogmento.com/")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://ui.cogmento.com/")
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:sumeet.desai@gmail.com")

# Adding a hyperlink re-stamps a brand-new (non-bordered) style, so restore
# the bordered hyperlink look (matching the pre-existing A3 / D7 cells) by
# re-pasting the format from D7, which already carries the correct bordered
# hyperlink style.
$ws.Range("D7").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. Column widths -- widen column C and give the new email/description
#    columns (G:H) a matching width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 23.67
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 23.67

# ---------------------------------------------------------------------------
# 6. Selection moves to F7 in the refreshed sheet.
# ---------------------------------------------------------------------------
$null = $ws.Range("F7").Select()
